$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Section caption: third method sampled from sailsim
$ws.Cells.Item(36, 1).Value = "dritte methode: gesampled von sailsim"

# --- Table 1 header (6kn wind table), typed in original column order ---
# (heading, heel degree, max speed kn, in in irons were typed first;
#  leeway was inserted later, matching the source shared-string order)
$ws.Cells.Item(41, 1).Value = "heading"
$ws.Cells.Item(41, 3).Value = "heel degree"
$ws.Cells.Item(41, 2).Value = "max speed kn"
$ws.Cells.Item(41, 5).Value = "in in irons"
$ws.Cells.Item(42, 5).Value = "yes"
$ws.Cells.Item(41, 4).Value = "leeway"

# Wind-condition captions for the three tables
$ws.Cells.Item(39, 1).Value = "6kn wind von 0 grad aus"
$ws.Cells.Item(60, 1).Value = "15kn wind  von 0 grad aus"
$ws.Cells.Item(76, 1).Value = "30kn wind  von 0 grad aus"

# Final header column added at the end
$ws.Cells.Item(41, 6).Value = "segel winkel geschätzt"

# --- Table 1 data (6kn wind von 0 grad aus) ---
$ws.Cells.Item(42, 1).Value = 360
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(42, 3).Value = 1
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(43, 1).Value = 345
$ws.Cells.Item(43, 2).Value = 0
$ws.Cells.Item(43, 3).Value = 1
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(43, 5).Value = "yes"
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(44, 1).Value = 340
$ws.Cells.Item(44, 2).Value = 0.9
$ws.Cells.Item(44, 3).Value = 1
$ws.Cells.Item(44, 4).Value = 4
$ws.Cells.Item(44, 6).Value = 3
$ws.Cells.Item(45, 1).Value = 320
$ws.Cells.Item(45, 2).Value = 2.1
$ws.Cells.Item(45, 3).Value = 2
$ws.Cells.Item(45, 4).Value = 4
$ws.Cells.Item(45, 6).Value = 10
$ws.Cells.Item(46, 1).Value = 308
$ws.Cells.Item(46, 2).Value = 2.6
$ws.Cells.Item(46, 3).Value = 2
$ws.Cells.Item(46, 4).Value = 4
$ws.Cells.Item(47, 1).Value = 290
$ws.Cells.Item(47, 2).Value = 3
$ws.Cells.Item(47, 3).Value = 2
$ws.Cells.Item(47, 4).Value = 3
$ws.Cells.Item(47, 6).Value = 25
$ws.Cells.Item(48, 1).Value = 284
$ws.Cells.Item(48, 2).Value = 2.9
$ws.Cells.Item(48, 3).Value = 2
$ws.Cells.Item(48, 4).Value = 3
$ws.Cells.Item(49, 1).Value = 270
$ws.Cells.Item(49, 2).Value = 2.9
$ws.Cells.Item(49, 3).Value = 1
$ws.Cells.Item(49, 4).Value = 3
$ws.Cells.Item(50, 1).Value = 250
$ws.Cells.Item(50, 2).Value = 3
$ws.Cells.Item(50, 3).Value = 1
$ws.Cells.Item(50, 4).Value = 3
$ws.Cells.Item(51, 1).Value = 238
$ws.Cells.Item(51, 2).Value = 2.3
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 2
$ws.Cells.Item(51, 6).Value = 45
$ws.Cells.Item(52, 1).Value = 222
$ws.Cells.Item(52, 2).Value = 2.2
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 4).Value = 2
$ws.Cells.Item(53, 1).Value = 208
$ws.Cells.Item(53, 2).Value = 1.9
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 1
$ws.Cells.Item(54, 1).Value = 193
$ws.Cells.Item(54, 2).Value = 1.8
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(54, 4).Value = 0
$ws.Cells.Item(55, 1).Value = 180
$ws.Cells.Item(55, 2).Value = 1.8
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(55, 4).Value = 0
$ws.Cells.Item(55, 6).Value = 90

# --- Table 2 header + data (15kn wind von 0 grad aus) ---
$ws.Cells.Item(62, 1).Value = "heading"
$ws.Cells.Item(62, 2).Value = "max speed kn"
$ws.Cells.Item(62, 3).Value = "heel degree"
$ws.Cells.Item(62, 4).Value = "leeway"
$ws.Cells.Item(62, 5).Value = "in in irons"
$ws.Cells.Item(62, 6).Value = "segel winkel geschätzt"
$ws.Cells.Item(63, 1).Value = 360
$ws.Cells.Item(63, 2).Value = 0
$ws.Cells.Item(63, 3).Value = 1
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = "yes"
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(64, 1).Value = 330
$ws.Cells.Item(64, 2).Value = 3.4
$ws.Cells.Item(64, 3).Value = 11
$ws.Cells.Item(64, 4).Value = 4
$ws.Cells.Item(64, 6).Value = 1
$ws.Cells.Item(65, 1).Value = 318
$ws.Cells.Item(65, 2).Value = 4.6
$ws.Cells.Item(65, 3).Value = 13
$ws.Cells.Item(65, 4).Value = 4
$ws.Cells.Item(65, 6).Value = 20
$ws.Cells.Item(66, 1).Value = 309
$ws.Cells.Item(66, 2).Value = 5.6
$ws.Cells.Item(66, 3).Value = 14
$ws.Cells.Item(66, 4).Value = 4
$ws.Cells.Item(66, 6).Value = 25
$ws.Cells.Item(67, 1).Value = 286
$ws.Cells.Item(67, 2).Value = 6.1
$ws.Cells.Item(67, 3).Value = 12
$ws.Cells.Item(67, 4).Value = 3
$ws.Cells.Item(67, 6).Value = 30
$ws.Cells.Item(68, 1).Value = 268
$ws.Cells.Item(68, 2).Value = 6.4
$ws.Cells.Item(68, 3).Value = 8
$ws.Cells.Item(68, 4).Value = 3
$ws.Cells.Item(68, 6).Value = 33
$ws.Cells.Item(69, 1).Value = 247
$ws.Cells.Item(69, 2).Value = 6.7
$ws.Cells.Item(69, 3).Value = 4
$ws.Cells.Item(69, 4).Value = 2
$ws.Cells.Item(69, 6).Value = 41
$ws.Cells.Item(70, 1).Value = 238
$ws.Cells.Item(70, 2).Value = 5.9
$ws.Cells.Item(70, 3).Value = 3
$ws.Cells.Item(70, 4).Value = 2
$ws.Cells.Item(70, 6).Value = 45
$ws.Cells.Item(71, 1).Value = 214
$ws.Cells.Item(71, 2).Value = 5.2
$ws.Cells.Item(71, 3).Value = 1
$ws.Cells.Item(71, 4).Value = 1
$ws.Cells.Item(71, 6).Value = 50
$ws.Cells.Item(72, 1).Value = 180
$ws.Cells.Item(72, 2).Value = 4.6
$ws.Cells.Item(72, 3).Value = 0
$ws.Cells.Item(72, 4).Value = 0
$ws.Cells.Item(72, 6).Value = 90

# --- Table 3 header + data (30kn wind von 0 grad aus) ---
$ws.Cells.Item(78, 1).Value = "heading"
$ws.Cells.Item(78, 2).Value = "max speed kn"
$ws.Cells.Item(78, 3).Value = "heel degree"
$ws.Cells.Item(78, 4).Value = "leeway"
$ws.Cells.Item(78, 5).Value = "in in irons"
$ws.Cells.Item(78, 6).Value = "segel winkel geschätzt"
$ws.Cells.Item(79, 1).Value = 360
$ws.Cells.Item(79, 2).Value = 0
$ws.Cells.Item(79, 3).Value = 1
$ws.Cells.Item(79, 4).Value = 0
$ws.Cells.Item(79, 5).Value = "yes"
$ws.Cells.Item(79, 6).Value = 0
$ws.Cells.Item(80, 1).Value = 329
$ws.Cells.Item(80, 2).Value = 3.8
$ws.Cells.Item(80, 3).Value = 38
$ws.Cells.Item(80, 4).Value = 4
$ws.Cells.Item(80, 6).Value = 5
$ws.Cells.Item(81, 1).Value = 316
$ws.Cells.Item(81, 2).Value = 3.7
$ws.Cells.Item(81, 3).Value = 40
$ws.Cells.Item(81, 4).Value = 3
$ws.Cells.Item(81, 6).Value = 20
$ws.Cells.Item(82, 1).Value = 300
$ws.Cells.Item(82, 2).Value = 3.9
$ws.Cells.Item(82, 3).Value = 40
$ws.Cells.Item(82, 4).Value = 3
$ws.Cells.Item(82, 6).Value = 25
$ws.Cells.Item(83, 1).Value = 279
$ws.Cells.Item(83, 2).Value = 4.5
$ws.Cells.Item(83, 3).Value = 32
$ws.Cells.Item(83, 4).Value = 3
$ws.Cells.Item(83, 6).Value = 35
$ws.Cells.Item(84, 1).Value = 262
$ws.Cells.Item(84, 2).Value = 5.4
$ws.Cells.Item(84, 3).Value = 21
$ws.Cells.Item(84, 4).Value = 2
$ws.Cells.Item(84, 6).Value = 40
$ws.Cells.Item(85, 1).Value = 250
$ws.Cells.Item(85, 2).Value = 6
$ws.Cells.Item(85, 3).Value = 13
$ws.Cells.Item(85, 4).Value = 2
$ws.Cells.Item(85, 6).Value = 45
$ws.Cells.Item(86, 1).Value = 220
$ws.Cells.Item(86, 2).Value = 6.3
$ws.Cells.Item(86, 3).Value = 3
$ws.Cells.Item(86, 4).Value = 1
$ws.Cells.Item(86, 6).Value = 50
$ws.Cells.Item(87, 1).Value = 200
$ws.Cells.Item(87, 2).Value = 6.1
$ws.Cells.Item(87, 3).Value = 1
$ws.Cells.Item(87, 4).Value = 1
$ws.Cells.Item(87, 6).Value = 80
$ws.Cells.Item(88, 1).Value = 180
$ws.Cells.Item(88, 2).Value = 6.1
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 6).Value = 90

# --- Column widths for the newly used columns B, C, E ---
$ws.Columns.Item(2).ColumnWidth = 11.59244791666667
$ws.Columns.Item(3).ColumnWidth = 12.30729166666667
$ws.Columns.Item(5).ColumnWidth = 9.16666666666667

# --- Print/page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection / active cell to match final saved view ---
$ws.Range("I48").Select()
